# Generate Report for Handback
# Swap the a.md / b.md rows' identities across the Overview / zh-cn / de-de
# sheets (file name + dependent status/date columns), update the hyperlink
# display text to match, and widen the "Status" columns that now hold the
# longer "not in sync" string.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "b.md"
$ws.Range("B2").Value = "e2e\b.md"
$ws.Range("G2").Value = "2017-02-21 04:05:08"

$ws.Range("A3").Value = "a.md"
$ws.Range("B3").Value = "e2e\a.md"
$ws.Range("G3").Value = "2017-02-21 04:05:08"

# Hyperlink display text follows the swapped file names (address/rId is
# left untouched, only the visible text changes).
$hyperlinks = @()
foreach ($hl in $ws.Hyperlinks) { $hyperlinks += $hl }
$hyperlinks[0].TextToDisplay = "e2e\b.md"
$hyperlinks[1].TextToDisplay = "e2e\a.md"

$ws.Columns.Item(5).ColumnWidth = 32.6
$ws.Columns.Item(6).ColumnWidth = 32.6

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "b.md"
$ws.Range("J2").Value = "b.md"
$ws.Range("C2").Value = "Handed back: not in sync with en-US"

$ws.Range("A3").Value = "a.md"
$ws.Range("J3").Value = "a.md"
$ws.Range("H3").Value = "2017-02-21 04:04:52"
$ws.Range("L3").Value = "2017-02-21 04:06:36"
$ws.Range("M3").Value = "TestHandback_201702211206"

$hyperlinks = @()
foreach ($hl in $ws.Hyperlinks) { $hyperlinks += $hl }
$hyperlinks[0].TextToDisplay = "b.md"
$hyperlinks[1].TextToDisplay = "b.md"
$hyperlinks[2].TextToDisplay = "a.md"
$hyperlinks[3].TextToDisplay = "a.md"

$ws.Columns.Item(3).ColumnWidth = 32.6

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "b.md"
$ws.Range("J2").Value = "b.md"
$ws.Range("C2").Value = "Handed back: not in sync with en-US"

$ws.Range("A3").Value = "a.md"
$ws.Range("J3").Value = "a.md"
$ws.Range("H3").Value = "2017-02-21 04:05:08"
$ws.Range("L3").Value = "2017-02-21 04:06:59"
$ws.Range("M3").Value = "TestHandback_201702211206"

$hyperlinks = @()
foreach ($hl in $ws.Hyperlinks) { $hyperlinks += $hl }
$hyperlinks[0].TextToDisplay = "b.md"
$hyperlinks[1].TextToDisplay = "b.md"
$hyperlinks[2].TextToDisplay = "a.md"
$hyperlinks[3].TextToDisplay = "a.md"

$ws.Columns.Item(3).ColumnWidth = 32.6
